$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Update the first few summary rows at the top of the table ---
$t.Rows.Item(1).Cells.Item(1).Range.Text = "0M"
$t.Rows.Item(2).Cells.Item(1).Range.Text = "0M"
$t.Rows.Item(3).Cells.Item(1).Range.Text = "0M"
$t.Rows.Item(4).Cells.Item(1).Range.Text = "218"
$t.Rows.Item(5).Cells.Item(1).Range.Text = "0.00002"

# Insert a new row right after row 5 and give it the new measurement.
$t.Rows.Add($t.Rows.Item(6)) | Out-Null
$t.Rows.Item(6).Cells.Item(1).Range.Text = "0.00007"

# The row that used to follow (now pushed to index 8, value "0.00003")
# is a stale duplicate and gets removed entirely.
$t.Rows.Item(8).Delete()

# Remaining stat rows that changed value.
$t.Rows.Item(9).Cells.Item(1).Range.Text = "0.00006"
$t.Rows.Item(10).Cells.Item(1).Range.Text = "0.00006"
$t.Rows.Item(11).Cells.Item(1).Range.Text = "0.00006"
$t.Rows.Item(12).Cells.Item(1).Range.Text = "0.00814"

# --- The final three rows previously packed 10 tab-separated values into
# one run each; collapse each back down to its single summary value. ---
$t.Rows.Item(44).Cells.Item(1).Range.Text = "100"
$t.Rows.Item(45).Cells.Item(1).Range.Text = "0.01"
$t.Rows.Item(46).Cells.Item(1).Range.Text = "367"
